$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# --- Cells changing between numeric and placeholder-text (`0` / `***.*`) representations ---
$ws.Range("G14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("H14").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("C15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Value = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("D20").Value = "'0"
$ws.Range("A20").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = "'***.*"
$ws.Range("A20").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("I22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

$ws.Range("I22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 3

$ws.Range("K22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -33.333333333333

$ws.Range("C26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Value = "'***.*"
$ws.Range("A26").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("A27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("I28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1

$ws.Range("I29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$ws.Range("D30").Value = "'0"
$ws.Range("A30").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = "'***.*"
$ws.Range("A30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -84
$ws.Range("L15").Value = 15
$ws.Range("M15").Value = -4.166666666666
$ws.Range("N15").Value = -63.492063492063
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 150
$ws.Range("J16").Value = 179
$ws.Range("K16").Value = -16.201117318435
$ws.Range("L16").Value = -10.714285714285
$ws.Range("M16").Value = -58.791208791208
$ws.Range("N16").Value = -91.263832265579
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -18.918918918918
$ws.Range("I17").Value = 375
$ws.Range("J17").Value = 350
$ws.Range("K17").Value = 7.142857142857
$ws.Range("L17").Value = 24.172185430463
$ws.Range("M17").Value = 4.166666666666
$ws.Range("N17").Value = -61.419753086419
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -37.931034482758
$ws.Range("I18").Value = 194
$ws.Range("J18").Value = 189
$ws.Range("K18").Value = 2.645502645502
$ws.Range("L18").Value = -36.393442622950
$ws.Range("M18").Value = -44.886363636363
$ws.Range("N18").Value = -81.663516068052
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -2.222222222222
$ws.Range("I19").Value = 345
$ws.Range("J19").Value = 318
$ws.Range("K19").Value = 8.490566037735
$ws.Range("L19").Value = 0.877192982456
$ws.Range("M19").Value = -2.542372881355
$ws.Range("N19").Value = -52.804377564979
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -62.5
$ws.Range("I20").Value = 80
$ws.Range("K20").Value = 2.564102564102
$ws.Range("L20").Value = 17.647058823529
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -82.417582417582
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 3.846153846153
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = -21.641791044776
$ws.Range("I21").Value = 1175
$ws.Range("J21").Value = 1134
$ws.Range("K21").Value = 3.615520282186
$ws.Range("L21").Value = -3.052805280528
$ws.Range("M21").Value = -24.968071519795
$ws.Range("N21").Value = -76.714229092350
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 5
$ws.Range("I22").Value = 24
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 4.347826086956
$ws.Range("L22").Value = 84.615384615384
$ws.Range("M22").Value = -17.241379310344
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 33.333333333333
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -18.75
$ws.Range("I23").Value = 198
$ws.Range("J23").Value = 206
$ws.Range("K23").Value = -3.883495145631
$ws.Range("L23").Value = 16.470588235294
$ws.Range("M23").Value = 7.027027027027
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -25.925925925925
$ws.Range("F24").Value = 156
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = 26.829268292682
$ws.Range("I24").Value = 1364
$ws.Range("J24").Value = 1070
$ws.Range("K24").Value = 27.476635514018
$ws.Range("L24").Value = 5.736434108527
$ws.Range("M24").Value = 61.803084223013
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = 2
$ws.Range("I25").Value = 563
$ws.Range("J25").Value = 419
$ws.Range("K25").Value = 34.367541766109
$ws.Range("L25").Value = 24.557522123893
$ws.Range("M25").Value = -32.655502392344
$ws.Range("L26").Value = 0
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -62.5
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -41.666666666666
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 30
$ws.Range("K28").Value = -41.176470588235
$ws.Range("L28").Value = -51.612903225806
$ws.Range("M28").Value = -40
$ws.Range("N28").Value = -87.704918032786
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 23
$ws.Range("K29").Value = -43.902439024390
$ws.Range("L29").Value = -53.061224489795
$ws.Range("M29").Value = -39.473684210526
$ws.Range("N29").Value = -89.201877934272
